# Update TestData_1 (sheet1) row 2: URL / UID / PWD columns
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TestData_1")
$ws2 = $wb.Worksheets.Item("TestData_2")

# --- TestData_1 data changes ---
$ws1.Range("E2").Value = "'http://admin.akshara.edu.pushvastech.in/login.php"
$ws1.Range("F2").Value = "'9885"
$ws1.Range("G2").Value = "'123"

# widen column E on TestData_1 to fit the new (longer) URL
$ws1.Range("E1").ColumnWidth = 44.5

# --- TestData_2 data changes ---
$ws2.Range("D2").Value = "N"

# --- selection / active sheet bookkeeping ---
$ws2.Activate()
$ws2.Range("D3").Select()

$ws1.Activate()
$ws1.Range("G3").Select()
